$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing header cell (H1) to the new header cells,
# then set their values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-20: I column = 1, J column = same value as H column
for ($r = 2; $r -le 20; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
